# Woo commerce to power office Documentation.docx -- update script
# Implements:
#  1. After "Documentation" title: add a line break + a new colored line
#     "open setting of plug in(wp-admin/admin.php?page=power_office_woocomerce)"
#  2. Step 3 paragraph: re-split the existing text run(s) into the same text
#     (no content change), matching the finer-grained run layout of the edit.
#  3. Step 4 paragraph: same re-split treatment for its text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: split the run(s) covering [start,end) away from their neighbours
# by toggling a formatting flag on and back off. This forces a run break
# at both ends of the sub-range while leaving the original formatting
# (including theme colors) completely untouched.
# ---------------------------------------------------------------------
function Split-Range([int]$start, [int]$end) {
    if ($end -gt $start) {
        $rr = $d.Range($start, $end)
        $rr.Font.Bold = $true
        $rr.Font.Bold = $false
    }
}

# =======================================================================
# CHANGE 1: title paragraph -- add line break + new URL-ish line
# =======================================================================

$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$insertPos = $titleRange.End - 1   # just before the paragraph mark

# --- line break, styled like the title (bold, green, 20pt) ---
$titleRange.InsertBreak(6)         # 6 = wdLineBreak
$afterBreak = $titleRange.End - 1
$brRange = $d.Range($insertPos, $afterBreak)
$brRange.Font.Bold = $true
$brRange.Font.Color = 5287936      # 00B050 (BGR-encoded)
$brRange.Font.Size = 20

$pos = $afterBreak

# --- the new blue 12pt text, built up run by run ---
$segments = @(
    "open setting of plug ",
    "in(",
    "wp",
    "-admin/",
    "admin.php?page",
    "=",
    "power_office_woocomerce",
    ")"
)

foreach ($seg in $segments) {
    $ip = $d.Range($pos, $pos)
    $ip.InsertAfter($seg)
    $segEnd = $pos + $seg.Length
    $segRange = $d.Range($pos, $segEnd)
    $segRange.Font.Color = 12611584  # 0070C0 (BGR-encoded)
    $segRange.Font.Size = 12
    $pos = $segEnd
}

# =======================================================================
# CHANGE 2: Step 3 paragraph -- re-split "upload plug-in zip file..." text
# =======================================================================

$full = $d.Content.Text
$anchor2 = "upload plug-in zip file activate plugin and go to plugin setting  Add power  office client and secret key make basic auth mention below"
$base2 = $full.IndexOf($anchor2)

$bounds2 = @(0, 33, 39, 50, 56, 57, 64, 69, 135)
for ($i = 0; $i -lt $bounds2.Length - 1; $i++) {
    $s2 = $base2 + $bounds2[$i]
    $e2 = $base2 + $bounds2[$i + 1]
    Split-Range $s2 $e2
}

# =======================================================================
# CHANGE 3: Step 4 paragraph -- re-split "add power office ..." text
# =======================================================================

$full = $d.Content.Text
$anchor3 = "add power office  api auth url and add woo commerce rest api key and power office key make basic authentication .we import all product orders and customers. when place order add product and create customer plugin automatically trigger this ."
$base3 = $full.IndexOf($anchor3)

$bounds3 = @(0, 10, 18, 21, 27, 30, 57, 60, 157, 161, 206, 212, 241)
for ($i = 0; $i -lt $bounds3.Length - 1; $i++) {
    $s3 = $base3 + $bounds3[$i]
    $e3 = $base3 + $bounds3[$i + 1]
    Split-Range $s3 $e3
}

Write-Output "edit complete"
